$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert a new column before column D for "Hours"
$ws.Range("D1").EntireColumn.Insert()

# 2. Header for the new column
$ws.Range("D2").Value = "Hours"

# 3. Fill in "Hours" values for existing rows (copy format from C column first, then
#    set number format + value). Row 3 gets the special bold/centered-row format
#    matching the rest of row 3 (font 5, border, right aligned) + date/time number
#    format (numFmtId 22), the rest copy the plain row style from C.
$ws.Range("C3").Copy()
$ws.Range("D3").PasteSpecial(-4122)
$ws.Range("D3").NumberFormat = "m/d/yy h:mm"
$ws.Range("D3").Value = "1 hour"

$ws.Range("C4").Copy()
$ws.Range("D4:D7").PasteSpecial(-4122)
$ws.Range("D4").Value = "2 hours 30"
$ws.Range("D5").Value = "4 hours 10"
$ws.Range("D6").Value = "4 hours"
$ws.Range("D7").Value = "3 hours 20"

# 4. New row 8 of data
$ws.Range("A8").Value = 41935
$ws.Range("B8").Value = "2:00am"
$ws.Range("C8").Value = "5:00am"
$ws.Range("D8").Value = "4 hours"
$ws.Range("E8").Value = 180
$ws.Range("F8").Value = "Programming"
$ws.Range("G8").Value = "Compelted basic mechancs of ball and player, added brick data."

# 5. Blank "Hours" cells for the remaining template rows (9-29) should carry the
#    plain row style (same as C9:C29) rather than the one inherited from the
#    column insert.
$ws.Range("C9").Copy()
$ws.Range("D9:D29").PasteSpecial(-4122)

# 6. Update the sheet view: scrolled to row 4, selection at G9
$ws.Range("G9").Select()
$excel.ActiveWindow.ScrollRow = 4
